$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record was added to the daily log. In the canonical sheet this
# shows up as a brand-new row inserted right before the former row 263,
# pushing every subsequent record down by one row (the former last row,
# 368, is duplicated as the new final row, 369).
$ws.Rows("263:263").Insert()

# Populate the newly inserted row with the new record's data. Columns that
# are identical to the record that used to occupy row 263 (now shifted to
# row 264) are simply carried over.
$ws.Range("A263").Value = 10
$ws.Range("B263").Value = "Vega Modelo de Temuco"
$ws.Range("C263").Value = "La Araucanía"
$ws.Range("D263").Value = 44524
$ws.Range("E263").Value = 9
$ws.Range("F263").Value = "Fruta"
$ws.Range("G263").Value = 100101
$ws.Range("H263").Value = "Berries"
$ws.Range("I263").Value = 100101007
$ws.Range("J263").Value = "Kiwi"
$ws.Range("K263").Value = "Hayward"
$ws.Range("L263").Value = "Especial"
$ws.Range("M263").Value = 95
$ws.Range("N263").Value = 24000
$ws.Range("O263").Value = 24000
$ws.Range("P263").Value = 24000
$ws.Range("Q263").Value = "`$/bandeja 18 kilos"
$ws.Range("R263").Value = "Región de O'Higgins"
$ws.Range("S263").Value = 1333
$ws.Range("T263").Value = 18
